$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 25
$ws.Range("B2").Value = 94
$ws.Range("C2").Value = 34513
$ws.Range("D2").Value = 90.00000000000001

# Row 3
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 71
$ws.Range("C3").Value = 34107

# Row 4
$ws.Range("A4").Value = 13
$ws.Range("B4").Value = 79
$ws.Range("C4").Value = 35119

# Row 5
$ws.Range("A5").Value = 22
$ws.Range("B5").Value = 99
$ws.Range("C5").Value = 39211

# Row 6
$ws.Range("A6").Value = 23
$ws.Range("B6").Value = 77
$ws.Range("C6").Value = 22825
